$d = $word.ActiveDocument
$p = $d.Paragraphs.Item(1)
$r = $p.Range
[void]$r.InsertXML('<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p w:rsidP="009168BC" w:rsidR="00E02A2B" w:rsidRDefault="00E02A2B"><w:pPr><w:tabs><w:tab w:pos="3119" w:val="left"/></w:tabs></w:pPr><w:r><w:t xml:space="preserve">Test link </w:t></w:r><w:r w:rsidR="00976419"><w:t>without</w:t></w:r><w:r><w:t xml:space="preserve"> bookmark : </w:t></w:r><w:r><w:rPr><w:b w:val="true"/><w:color w:val="FF0000"/></w:rPr><w:t>dangling reference for bookmark bookmark1</w:t></w:r><w:r/><w:r w:rsidR="32E4AE5A86ABB53458D48AF6A0222314"><w:fldChar w:fldCharType="begin"/></w:r><w:r w:rsidR="32E4AE5A86ABB53458D48AF6A0222314"><w:instrText xml:space="preserve"/></w:r><w:r w:rsidR="32E4AE5A86ABB53458D48AF6A0222314"><w:fldChar w:fldCharType="separate"/></w:r><w:r w:rsidR="32E4AE5A86ABB53458D48AF6A0222314"><w:rPr><w:b w:val="true"/><w:noProof/></w:rPr><w:t>a reference to bookmark1</w:t></w:r><w:r w:rsidR="32E4AE5A86ABB53458D48AF6A0222314"><w:fldChar w:fldCharType="end"/></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>')
